$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text edits on the existing agenda grid -------------------------------
$ws.Range("C2").Value = "Caleb, Mitch, Corinne, Mike"
$ws.Range("F15").Value = "Caleb"
$ws.Range("D18").Value = "Investigate LID alternatives for new dev'pt (ANC3)"

# --- Remove the "Bonus" block (rows 26-36), keeping row 25's heading -------
$ws.Rows("26:36").Delete() | Out-Null

# Replace row 25's plain heading with the two-tone rich text label
$ws.Range("B25").Value = "Bonus Lectures: Model Calibration, Water Quality"
$prefixLen = ("Bonus Lectures").Length
$fullLen = ("Bonus Lectures: Model Calibration, Water Quality").Length
$ws.Range("B25").Characters($prefixLen + 1, $fullLen - $prefixLen).Font.Bold = $false

# --- Column width / autosize tweaks ----------------------------------------
$ws.Columns("D").ColumnWidth = 40.42
$ws.Columns("F").ColumnWidth = 12.59

# --- View tweaks -------------------------------------------------------------
$excel.ActiveWindow.Zoom = 100
